# Scheduled runner update: refresh market price / profit figures on each class sheet
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 3780.5
$ws.Range("J29").Value = 6167.5
$ws.Range("L29").Value = 18502.5
$ws.Range("N29").Value = -19064.5
$ws.Range("H58").Value = 2071.5
$ws.Range("I58").Value = 107.25
$ws.Range("J58").Value = 6000
$ws.Range("K58").Value = 321.75
$ws.Range("L58").Value = 18000
$ws.Range("M58").Value = -171.75
$ws.Range("N58").Value = -18300
$ws.Range("H61").Value = 157.5
$ws.Range("I61").Value = 157.5
$ws.Range("K61").Value = 472.5
$ws.Range("M61").Value = -300.5
$ws.Range("H70").Value = 46185.625
$ws.Range("I70").Value = 2295.5454
$ws.Range("J70").Value = 142743.8
$ws.Range("K70").Value = 6886.6362
$ws.Range("L70").Value = 428231.4
$ws.Range("M70").Value = -6616.6362
$ws.Range("N70").Value = -428771.4
$ws.Range("H73").Value = 46185.625
$ws.Range("I73").Value = 2295.5454
$ws.Range("J73").Value = 142743.8
$ws.Range("K73").Value = 6886.6362
$ws.Range("L73").Value = 428231.4
$ws.Range("M73").Value = -5950.6362
$ws.Range("N73").Value = -430103.4
$ws.Range("H87").Value = 69850
$ws.Range("J87").Value = 69850
$ws.Range("L87").Value = 69850
$ws.Range("N87").Value = -72346
$ws.Range("H90").Value = 69850
$ws.Range("J90").Value = 69850
$ws.Range("L90").Value = 209550
$ws.Range("N90").Value = -222030
$ws.Range("H138").Value = 9247.427
$ws.Range("J138").Value = 8678.796
$ws.Range("L138").Value = 26036.388
$ws.Range("N138").Value = -36316.388

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3216.6924
$ws.Range("I2").Value = 2531.7
$ws.Range("K2").Value = 2531.7
$ws.Range("M2").Value = -2418.7
$ws.Range("H32").Value = 18835.314
$ws.Range("J32").Value = 30664.75
$ws.Range("L32").Value = 30664.75
$ws.Range("N32").Value = -31238.75
$ws.Range("H45").Value = 3045.647
$ws.Range("I45").Value = 2677.6
$ws.Range("K45").Value = 2677.6
$ws.Range("M45").Value = -2300.6
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("N62").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("N65").Value = 0
$ws.Range("H116").Value = 3216.6924
$ws.Range("I116").Value = 2531.7
$ws.Range("K116").Value = 2531.7
$ws.Range("M116").Value = -237.6999999999998
$ws.Range("H139").Value = 50650
$ws.Range("I139").Value = 50650
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 50650
$ws.Range("L139").ClearContents()
$ws.Range("M139").Value = -45510
$ws.Range("N139").Value = 0

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3216.6924
$ws.Range("I3").Value = 2531.7
$ws.Range("K3").Value = 2531.7
$ws.Range("M3").Value = -2417.7
$ws.Range("H20").Value = 2395.0417
$ws.Range("J20").Value = 3278.75
$ws.Range("L20").Value = 3278.75
$ws.Range("N20").Value = -3772.75
$ws.Range("H64").Value = 2206.6155
$ws.Range("I64").Value = 2274.125
$ws.Range("J64").Value = 2098.6
$ws.Range("K64").Value = 2274.125
$ws.Range("L64").Value = 2098.6
$ws.Range("M64").Value = -2049.125
$ws.Range("N64").Value = -2548.6
$ws.Range("H67").Value = 2206.6155
$ws.Range("I67").Value = 2274.125
$ws.Range("J67").Value = 2098.6
$ws.Range("K67").Value = 2274.125
$ws.Range("L67").Value = 2098.6
$ws.Range("M67").Value = -1494.125
$ws.Range("N67").Value = -3658.6
$ws.Range("H86").Value = 6087.5
$ws.Range("I86").Value = 5971.4287
$ws.Range("J86").Value = 6900
$ws.Range("K86").Value = 5971.4287
$ws.Range("L86").Value = 6900
$ws.Range("M86").Value = -4848.4287
$ws.Range("N86").Value = -9146
$ws.Range("H89").Value = 6087.5
$ws.Range("I89").Value = 5971.4287
$ws.Range("J89").Value = 6900
$ws.Range("K89").Value = 29857.1435
$ws.Range("L89").Value = 34500
$ws.Range("M89").Value = -24241.1435
$ws.Range("N89").Value = -45732

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5968.05
$ws.Range("I31").Value = 4882.9
$ws.Range("J31").Value = 7053.2
$ws.Range("K31").Value = 4882.9
$ws.Range("L31").Value = 7053.2
$ws.Range("M31").Value = -4587.9
$ws.Range("N31").Value = -7643.2
$ws.Range("H34").Value = 5968.05
$ws.Range("I34").Value = 4882.9
$ws.Range("J34").Value = 7053.2
$ws.Range("K34").Value = 4882.9
$ws.Range("L34").Value = 7053.2
$ws.Range("M34").Value = -4680.9
$ws.Range("N34").Value = -7457.2
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").ClearContents()
$ws.Range("N50").Value = 0
$ws.Range("H60").Value = 43999
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").ClearContents()
$ws.Range("N74").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").ClearContents()
$ws.Range("N77").Value = 0
$ws.Range("H132").Value = 2321.8057
$ws.Range("I132").Value = 2166.7585
$ws.Range("K132").Value = 6500.2755
$ws.Range("M132").Value = -3970.2755

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1777.5
$ws.Range("J131").Value = 1996.4286
$ws.Range("L131").Value = 5989.2858
$ws.Range("N131").Value = -16069.2858

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").ClearContents()
$ws.Range("N59").Value = 0
$ws.Range("H80").Value = 9998.5
$ws.Range("J80").Value = 9998.5
$ws.Range("L80").Value = 9998.5
$ws.Range("N80").Value = -11994.5
$ws.Range("H83").Value = 9998.5
$ws.Range("J83").Value = 9998.5
$ws.Range("L83").Value = 49992.5
$ws.Range("N83").Value = -59976.5

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5027968
$ws.Range("J2").Value = 50935.5
$ws.Range("L2").Value = 50935.5
$ws.Range("N2").Value = -51159.5
$ws.Range("H7").Value = 4999.25
$ws.Range("J7").Value = 4999
$ws.Range("L7").Value = 4999
$ws.Range("N7").Value = -5223
$ws.Range("H46").Value = 2944.389
$ws.Range("J46").Value = 2999.9167
$ws.Range("L46").Value = 2999.9167
$ws.Range("N46").Value = -3375.9167
$ws.Range("H82").Value = 2284.077
$ws.Range("J82").Value = 2133.3333
$ws.Range("L82").Value = 2133.3333
$ws.Range("N82").Value = -2855.3333
$ws.Range("H85").Value = 2284.077
$ws.Range("J85").Value = 2133.3333
$ws.Range("L85").Value = 2133.3333
$ws.Range("N85").Value = -4629.3333
$ws.Range("H99").Value = 13839
$ws.Range("I99").Value = 13839
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 13839
$ws.Range("L99").ClearContents()
$ws.Range("M99").Value = -10844
$ws.Range("N99").Value = 0
$ws.Range("H126").Value = 4999.25
$ws.Range("J126").Value = 4999
$ws.Range("L126").Value = 14997
$ws.Range("N126").Value = -19937
$ws.Range("H132").Value = 5208
$ws.Range("I132").Value = 2929.3333
$ws.Range("K132").Value = 8787.999899999999
$ws.Range("M132").Value = -6257.999899999999
$ws.Range("H136").Value = 2100.4167
$ws.Range("I136").Value = 2190.6
$ws.Range("K136").Value = 6571.799999999999
$ws.Range("M136").Value = -4021.799999999999

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 74170.64
$ws.Range("I14").Value = 127200.5
$ws.Range("J14").Value = 3464.1667
$ws.Range("K14").Value = 127200.5
$ws.Range("L14").Value = 3464.1667
$ws.Range("M14").Value = -127032.5
$ws.Range("N14").Value = -3800.1667
$ws.Range("H109").Value = 39999
$ws.Range("J109").Value = 39999
$ws.Range("L109").Value = 39999
$ws.Range("N109").Value = -42773
$ws.Range("H115").Value = 125000
$ws.Range("J115").Value = 125000
$ws.Range("L115").Value = 125000
$ws.Range("N115").Value = -128134
$ws.Range("H136").Value = 57703.777
$ws.Range("I136").Value = 1369.9286
$ws.Range("K136").Value = 4109.7858
$ws.Range("M136").Value = -1559.7858
